# Apply cryptocurrency price/volume updates from the GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.024.03'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '2.274.55'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '''231.12'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('D6').Value = '''0.631'
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('D7').Value = '''63.59'
$ws.Range('E7').Value = '  +2.27%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '''0.448'
$ws.Range('E9').Value = '  +7.03%  '
$ws.Range('E10').Value = '  +7.73%  '
$ws.Range('D11').Value = '''57.58'
$ws.Range('E11').Value = '  -1.04%  '
$ws.Range('D12').Value = '''27.35'
$ws.Range('E12').Value = '  +13.19%  '
$ws.Range('E13').Value = '  +1.59%  '
$ws.Range('D14').Value = '2.613.77'
$ws.Range('E14').Value = '  -0.92%  '
$ws.Range('D15').Value = '''15.77'
$ws.Range('E15').Value = '  -0.64%  '
$ws.Range('E16').Value = '  +6.40%  '
$ws.Range('D17').Value = '''0.836'
$ws.Range('E17').Value = '  +2.45%  '
$ws.Range('D18').Value = '2.279.73'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('D19').Value = '43.902.08'
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('E20').Value = '  +7.73%  '
$ws.Range('D21').Value = '''73.77'
$ws.Range('E21').Value = '  +0.43%  '
$ws.Range('E22').Value = '  -2.42%  '
$ws.Range('D23').Value = '''252.73'
$ws.Range('E23').Value = '  +0.90%  '
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('E25').Value = '  -4.34%  '
$ws.Range('E26').Value = '  -4.72%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '''10.05'
$ws.Range('E27').Value = '  +1.50%  '
$ws.Range('B28').Value = 'WEMIXToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D28').Value = '''3.31'
$ws.Range('E28').Value = '  +24.46%  '
$ws.Range('D29').Value = '''171.66'
$ws.Range('E29').Value = '  +1.09%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = '''0.139'
$ws.Range('E30').Value = '  -2.37%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '''20.91'
$ws.Range('E31').Value = '  +1.24%  '
$ws.Range('D32').Value = '''1.39'
$ws.Range('E32').Value = '  -6.32%  '
$ws.Range('E33').Value = '  +2.62%  '
$ws.Range('D34').Value = '''0.0700'
$ws.Range('E34').Value = '  +5.96%  '
$ws.Range('D35').Value = '''4.80'
$ws.Range('E35').Value = '  +1.10%  '
$ws.Range('D36').Value = '''4.87'
$ws.Range('E36').Value = '  -2.92%  '
$ws.Range('D37').Value = '''3.81'
$ws.Range('E37').Value = '  +3.81%  '
$ws.Range('D38').Value = '''6.54'
$ws.Range('E38').Value = '  -0.27%  '
$ws.Range('D39').Value = '''2.32'
$ws.Range('E39').Value = '  -5.50%  '
$ws.Range('E40').Value = '  +2.59%  '
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').Value = '''0.000230'
$ws.Range('E42').Value = '  +5.92%  '
$ws.Range('D43').Value = '''17.64'
$ws.Range('E43').Value = '  +5.35%  '
$ws.Range('D44').Value = '''0.0983'
$ws.Range('E44').Value = '  +0.86%  '
$ws.Range('E45').Value = '  -5.79%  '
$ws.Range('D46').Value = '''10.46'
$ws.Range('E46').Value = '  +11.24%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '''98.26'
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').Value = '''1.21'
$ws.Range('E48').Value = '  -0.78%  '
$ws.Range('D49').Value = '''4.39'
$ws.Range('E49').Value = '  -4.56%  '
$ws.Range('D50').Value = '1.449.51'
$ws.Range('E50').Value = '  -1.73%  '
$ws.Range('D51').Value = '''2.31'
$ws.Range('E51').Value = '  +0.46%  '
